$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old I:K columns that are no longer used (values + formatting)
$ws.Range("I1:K4").Clear()

# A1 was always empty
$ws.Range("A1").ClearContents()

# Header row (now numeric, not shared-string labels), columns B..H
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 1
$ws.Range("D1").Value = 0
$ws.Range("E1").Value = 1
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = 3
$ws.Range("H1").Value = 4

# Row 2 - FE
$ws.Range("A2").Value = "FE"
$ws.Range("B2").Value = 0.83
$ws.Range("C2").Value = 1.85
$ws.Range("D2").Value = 0.89
$ws.Range("E2").Value = 0.02
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 0.9
$ws.Range("H2").Value = 0.66

# Row 3 - FE+Disg
$ws.Range("A3").Value = "FE+Disg"
$ws.Range("B3").Value = 0.83
$ws.Range("C3").Value = 1.85
$ws.Range("D3").Value = 0.61
$ws.Range("E3").Value = 0.42
$ws.Range("F3").Value = 0.16
$ws.Range("G3").Value = 0.93
$ws.Range("H3").Value = 0.5600000000000001

# Row 4 - FE+Disg+Var
$ws.Range("A4").Value = "FE+Disg+Var"
$ws.Range("B4").Value = 0.83
$ws.Range("C4").Value = 1.85
$ws.Range("D4").Value = 0.58
$ws.Range("E4").Value = 0.5600000000000001
$ws.Range("F4").Value = 0.45
$ws.Range("G4").Value = 0.91
$ws.Range("H4").Value = 0.11
